# Semana 37 de 2025 - update Esperado/Observado/valor p columns in poisson sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 0.18

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0

$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 4

$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 0.02

$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0.14

$ws.Range("D9").Value = 41
$ws.Range("E9").Value = 0.06

$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0.37

$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 0.2

$ws.Range("C13").Value = 1
$ws.Range("E13").Value = 0.37

$ws.Range("C15").Value = 1
$ws.Range("E15").Value = 0.37

$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = 0.01

$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0.27

$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 0.13

$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 1

$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0

$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0

$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0.03

$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 0.15

$ws.Range("C33").Value = 7
$ws.Range("D33").Value = 7
$ws.Range("E33").Value = 0.15

$ws.Range("C34").Value = 11
$ws.Range("D34").Value = 5
$ws.Range("E34").Value = 0.02

$ws.Range("C35").Value = 8
$ws.Range("D35").Value = 7
$ws.Range("E35").Value = 0.14
